$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Trening" header in column F, matching the style of the existing header row
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Trening"

# Register the date/time number format (first lowercase, then the real uppercase
# code actually used) on A2 before filling the rest of the timestamp column so the
# style table ends up with both format strings declared.
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A3:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(2, 1).Value = 45684.59234166667
$ws.Cells.Item(2, 2).Value = 578.3
$ws.Cells.Item(2, 3).Value = 10.89
$ws.Cells.Item(2, 4).Value = 1.64991353239332
$ws.Cells.Item(2, 5).Value = "10-15"
$ws.Cells.Item(2, 6).Value = "Duża Gra"

$ws.Cells.Item(3, 1).Value = 45684.5925650463
$ws.Cells.Item(3, 2).Value = 597.6
$ws.Cells.Item(3, 3).Value = 10.04
$ws.Cells.Item(3, 4).Value = 2.001077464648655
$ws.Cells.Item(3, 5).Value = "10-15"
$ws.Cells.Item(3, 6).Value = "Duża Gra"

$ws.Cells.Item(4, 1).Value = 45684.59322592593
$ws.Cells.Item(4, 2).Value = 654.7
$ws.Cells.Item(4, 3).Value = 10.76
$ws.Cells.Item(4, 4).Value = 1.665620037487576
$ws.Cells.Item(4, 5).Value = "10-15"
$ws.Cells.Item(4, 6).Value = "Duża Gra"

$ws.Cells.Item(5, 1).Value = 45684.59073518519
$ws.Cells.Item(5, 2).Value = 439.5
$ws.Cells.Item(5, 3).Value = 7.3
$ws.Cells.Item(5, 4).Value = 1.600568413734434
$ws.Cells.Item(5, 5).Value = "5-10"
$ws.Cells.Item(5, 6).Value = "Duża Gra"

$ws.Cells.Item(6, 1).Value = 45684.59093888889
$ws.Cells.Item(6, 2).Value = 457.1
$ws.Cells.Item(6, 3).Value = 5.49
$ws.Cells.Item(6, 4).Value = 1.652362687247141
$ws.Cells.Item(6, 5).Value = "5-10"
$ws.Cells.Item(6, 6).Value = "Duża Gra"

$ws.Cells.Item(7, 1).Value = 45684.59256273148
$ws.Cells.Item(7, 2).Value = 597.4
$ws.Cells.Item(7, 3).Value = 8.75
$ws.Cells.Item(7, 4).Value = 2.043854781559534
$ws.Cells.Item(7, 5).Value = "5-10"
$ws.Cells.Item(7, 6).Value = "Duża Gra"

$ws.Cells.Item(8, 1).Value = 45684.60080462963
$ws.Cells.Item(8, 2).Value = 1309.5
$ws.Cells.Item(8, 3).Value = 14.89
$ws.Cells.Item(8, 4).Value = 3.255911588668824
$ws.Cells.Item(8, 5).Value = "10-15"
$ws.Cells.Item(8, 6).Value = "Mała Gra"

$ws.Cells.Item(9, 1).Value = 45684.60174675926
$ws.Cells.Item(9, 2).Value = 1390.9
$ws.Cells.Item(9, 3).Value = 13.5
$ws.Cells.Item(9, 4).Value = 3.104380096708025
$ws.Cells.Item(9, 5).Value = "10-15"
$ws.Cells.Item(9, 6).Value = "Mała Gra"

$ws.Cells.Item(10, 1).Value = 45684.60473055555
$ws.Cells.Item(10, 2).Value = 1648.7
$ws.Cells.Item(10, 3).Value = 14.34
$ws.Cells.Item(10, 4).Value = 3.56577604157584
$ws.Cells.Item(10, 5).Value = "10-15"
$ws.Cells.Item(10, 6).Value = "Mała Gra"

$ws.Cells.Item(11, 1).Value = 45684.59967962963
$ws.Cells.Item(11, 2).Value = 1212.3
$ws.Cells.Item(11, 3).Value = 8.98
$ws.Cells.Item(11, 4).Value = 2.569038936070033
$ws.Cells.Item(11, 5).Value = "5-10"
$ws.Cells.Item(11, 6).Value = "Mała Gra"

$ws.Cells.Item(12, 1).Value = 45684.6008
$ws.Cells.Item(12, 2).Value = 1309.1
$ws.Cells.Item(12, 3).Value = 9.87
$ws.Cells.Item(12, 4).Value = 2.441896906920842
$ws.Cells.Item(12, 5).Value = "5-10"
$ws.Cells.Item(12, 6).Value = "Mała Gra"

$ws.Cells.Item(13, 1).Value = 45684.60174212963
$ws.Cells.Item(13, 2).Value = 1390.5
$ws.Cells.Item(13, 3).Value = 8.97
$ws.Cells.Item(13, 4).Value = 2.58960109097617
$ws.Cells.Item(13, 5).Value = "5-10"
$ws.Cells.Item(13, 6).Value = "Mała Gra"
